$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "308.40"
Set-TextValue "E2" "0.17%"
Set-TextValue "D3" "40.86"
Set-TextValue "E3" "3.28%"
Set-TextValue "D4" "5.116"
Set-TextValue "E4" "0.38%"
Set-TextValue "D5" "0.07623"
Set-TextValue "E5" "-0.84%"
Set-TextValue "B6" "FTXToken"
Set-TextValue "C6" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D6" "1.601"
Set-TextValue "E6" "-0.27%"
Set-TextValue "B7" "BTSEToken"
Set-TextValue "C7" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D7" "2.460"
Set-TextValue "E7" "1.67%"
Set-TextValue "B8" "MXToken"
Set-TextValue "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.9021"
Set-TextValue "E8" "0.30%"
Set-TextValue "B9" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C9" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D9" "0.1122"
Set-TextValue "E9" "11.69%"
Set-TextValue "B10" "WazirX"
Set-TextValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1797"
Set-TextValue "E10" "3.64%"
Set-TextValue "B11" "MandalaExchangeToken"
Set-TextValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.09192"
Set-TextValue "E11" "2.26%"
Set-TextValue "B12" "BitrueCoin"
Set-TextValue "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.04158"
Set-TextValue "E12" "-5.78%"
Set-TextValue "B13" "BitMartToken"
Set-TextValue "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.1053"
Set-TextValue "E13" "-0.01%"
Set-TextValue "B14" "BitForexToken"
Set-TextValue "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001264"
Set-TextValue "E14" "0.51%"
Set-TextValue "B15" "TigerCash"
Set-TextValue "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D15" "0.005894"
Set-TextValue "E15" "-0.18%"
Set-TextValue "B16" "LEO"
Set-TextValue "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D16" "3.342"
Set-TextValue "E16" "-0.35%"
Set-TextValue "B17" "GateToken"
Set-TextValue "C17" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D17" "4.247"
Set-TextValue "E17" "0.34%"
Set-TextValue "D18" "0.3315"
Set-TextValue "E18" "0.00%"
Set-TextValue "D19" "6.639"
Set-TextValue "E19" "-6.08%"
Set-TextValue "E20" "1.16%"
Set-TextValue "D22" "0.04075"
Set-TextValue "E22" "-1.75%"
Set-TextValue "D23" "0.001247"
Set-TextValue "E23" "3.06%"
Set-TextValue "D24" "0.004106"
Set-TextValue "E24" "1.20%"
Set-TextValue "E25" "-0.23%"
Set-TextValue "D26" "0.0003747"
Set-TextValue "D38" "0.02409"
Set-TextValue "E38" "3.08%"
Set-TextValue "D39" "0.05197"
Set-TextValue "E39" "0.72%"
Set-TextValue "D40" "0.007781"
Set-TextValue "E40" "-1.94%"
Set-TextValue "E41" "-1.57%"
Set-TextValue "E42" "11.93%"
Set-TextValue "E43" "-0.20%"
Set-TextValue "D44" "0.007721"
Set-TextValue "E44" "-6.26%"
Set-TextValue "E45" "-7.71%"
Set-TextValue "D46" "0.00006968"
Set-TextValue "E46" "6.83%"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "-0.25%"
Set-TextValue "D48" "0.04170"
Set-TextValue "E48" "1,103.33%"
Set-TextValue "D50" "0.00002101"
Set-TextValue "E50" "-0.25%"
Set-TextValue "D51" "0.0002001"
Set-TextValue "E51" "-0.25%"
